$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 3.28683254486478793
$ws.Range("C2").Value2 = 1.65577808226027101
$ws.Range("D2").Value2 = 0.1494219747398047
$ws.Range("E2").Value2 = 0.49423653606076973
$ws.Range("G2").Value2 = 5.58626913792563418

$ws.Range("B3").Value2 = 3.28683254486478793
$ws.Range("C3").Value2 = 1.65577808226027101
$ws.Range("D3").Value2 = 3.53776164880671917
$ws.Range("E3").Value2 = 0.49423653606076973
$ws.Range("G3").Value2 = 8.97460881199254779

$ws.Range("B4").Value2 = 3.28683254486478793
$ws.Range("C4").Value2 = 1.65577808226027101
$ws.Range("D4").Value2 = 0.1494219747398047
$ws.Range("E4").Value2 = 0.49423653606076973
$ws.Range("G4").Value2 = 5.58626913792563418

$ws.Range("B5").Value2 = 3.28683254486478793
$ws.Range("C5").Value2 = 1.65577808226027101
$ws.Range("D5").Value2 = 0.1494219747398047
$ws.Range("E5").Value2 = 0.49423653606076973
$ws.Range("G5").Value2 = 5.58626913792563418

$ws.Range("B6").Value2 = 0.000001295275857016165
$ws.Range("C6").Value2 = 0.002571899574220771
$ws.Range("D6").Value2 = 3.53776164880671917
$ws.Range("E6").Value2 = 0.49423653606076973
$ws.Range("G6").Value2 = 4.03457137971756641

$ws.Range("B7").Value2 = 1.4553620445145421
$ws.Range("C7").Value2 = 1.65577808226027101
$ws.Range("D7").Value2 = 0.1494219747398047
$ws.Range("E7").Value2 = 0.49423653606076973
$ws.Range("G7").Value2 = 3.75479863757538679

$ws.Range("B8").Value2 = 0.66065244103595555
$ws.Range("C8").Value2 = 1.65577808226027101
$ws.Range("D8").Value2 = 0.1494219747398047
$ws.Range("E8").Value2 = 0.49423653606076973
$ws.Range("G8").Value2 = 2.96008903409680091

$ws.Range("B9").Value2 = 0.1190320826869504
$ws.Range("C9").Value2 = 0.30682122725969802
$ws.Range("D9").Value2 = 0.75274326777386413
$ws.Range("E9").Value2 = 10.1924530069365602
$ws.Range("G9").Value2 = 11.3710495846570705

$ws.Range("B10").Value2 = 3.28683254486478793
$ws.Range("C10").Value2 = 1.65577808226027101
$ws.Range("D10").Value2 = 0.75274326777386413
$ws.Range("E10").Value2 = 0.49423653606076973
$ws.Range("G10").Value2 = 6.18959043095969363

$ws.Range("B11").Value2 = 0.29177164025654623
$ws.Range("C11").Value2 = 1.65577808226027101
$ws.Range("D11").Value2 = 0.1494219747398047
$ws.Range("E11").Value2 = 10.1924530069365602
$ws.Range("G11").Value2 = 12.28942470419318056

$ws.Range("B12").Value2 = 1.4553620445145421
$ws.Range("C12").Value2 = 1.65577808226027101
$ws.Range("D12").Value2 = 0.1494219747398047
$ws.Range("E12").Value2 = 0.49423653606076973
$ws.Range("G12").Value2 = 3.75479863757538679

$ws.Range("B13").Value2 = 1.4553620445145421
$ws.Range("C13").Value2 = 1.65577808226027101
$ws.Range("D13").Value2 = 0.75274326777386413
$ws.Range("E13").Value2 = 0.49423653606076973
$ws.Range("G13").Value2 = 4.35811993060944669

$ws.Range("B14").Value2 = 0.29177164025654623
$ws.Range("C14").Value2 = 0.30682122725969802
$ws.Range("D14").Value2 = 0.1494219747398047
$ws.Range("E14").Value2 = 0.49423653606076973
$ws.Range("G14").Value2 = 1.24225137831681898

$ws.Range("B15").Value2 = 1.4553620445145421
$ws.Range("C15").Value2 = 1.65577808226027101
$ws.Range("D15").Value2 = 0.75274326777386413
$ws.Range("E15").Value2 = 0.49423653606076973
$ws.Range("G15").Value2 = 4.35811993060944669

$ws.Range("B16").Value2 = 0.04271373187048222
$ws.Range("C16").Value2 = 0.04071648406533734
$ws.Range("D16").Value2 = 0.1494219747398047
$ws.Range("E16").Value2 = 0.49423653606076973
$ws.Range("G16").Value2 = 0.72708872673639391

$ws.Range("B17").Value2 = 1.4553620445145421
$ws.Range("C17").Value2 = 1.65577808226027101
$ws.Range("D17").Value2 = 3.53776164880671917
$ws.Range("E17").Value2 = 0.49423653606076973
$ws.Range("G17").Value2 = 7.14313831164230173
